# Update "想去人数" (F) and "最低票价" (G) figures on the "展览" and
# "全部类型" worksheets to reflect newly generated output (gh-pages
# rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet ---
$ws1.Range("F3").Value = 560
$ws1.Range("G3").Value = 99
$ws1.Range("F4").Value = 1129
$ws1.Range("F5").Value = 118
$ws1.Range("F6").Value = 69
$ws1.Range("F8").Value = 57
$ws1.Range("F9").Value = 1157
$ws1.Range("F10").Value = 16278
$ws1.Range("F11").Value = 276
$ws1.Range("F12").Value = 201
$ws1.Range("F14").Value = 6350
$ws1.Range("F17").Value = 78
$ws1.Range("F18").Value = 20
$ws1.Range("F21").Value = 29
$ws1.Range("F24").Value = 30
$ws1.Range("F26").Value = 15
$ws1.Range("F27").Value = 219
$ws1.Range("F28").Value = 889
$ws1.Range("F29").Value = 52
$ws1.Range("F30").Value = 5046
$ws1.Range("F31").Value = 496
$ws1.Range("F32").Value = 11303
$ws1.Range("F33").Value = 1244
$ws1.Range("F34").Value = 18
$ws1.Range("F36").Value = 200
$ws1.Range("F37").Value = 3834
$ws1.Range("F38").Value = 269

# --- 全部类型 sheet ---
$ws4.Range("F3").Value = 560
$ws4.Range("G3").Value = 99
$ws4.Range("F4").Value = 1129
$ws4.Range("F5").Value = 118
$ws4.Range("F6").Value = 69
$ws4.Range("F8").Value = 57
$ws4.Range("F9").Value = 1157
$ws4.Range("F10").Value = 16278
$ws4.Range("F11").Value = 276
$ws4.Range("F12").Value = 201
$ws4.Range("F14").Value = 6350
$ws4.Range("F17").Value = 78
$ws4.Range("F18").Value = 20
$ws4.Range("F21").Value = 29
$ws4.Range("F24").Value = 30
$ws4.Range("F26").Value = 15
$ws4.Range("F27").Value = 219
$ws4.Range("F28").Value = 889
$ws4.Range("F29").Value = 52
$ws4.Range("F30").Value = 5046
$ws4.Range("F31").Value = 496
$ws4.Range("F33").Value = 11303
$ws4.Range("F34").Value = 1244
$ws4.Range("F35").Value = 18
$ws4.Range("F37").Value = 200
$ws4.Range("F38").Value = 3834
$ws4.Range("F39").Value = 269

$wb.Save()
